$wb = $excel.ActiveWorkbook

# --- Deselect the currently active sheets / move selection (cosmetic, matches target) ---
$wsHealth = $wb.Worksheets.Item(1)
$wsHealth.Range("I7").Select() | Out-Null

$wsCalib = $wb.Worksheets.Item(2)
$wsCalib.Activate()
$wsCalib.Range("I7").Select() | Out-Null

# --- Add the new "Final results" sheet after "Calibration" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Final results"
$ws.Activate()

# --- Column width ---
$ws.Columns.Item(4).ColumnWidth = 11.27

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 29
$ws.Rows.Item(2).RowHeight = 43.5

# --- Header row 1 ---
$ws.Range("C1").Value = "LIT RATES"
$ws.Range("D1").Value = "NO FIRE"
$ws.Range("E1").Value = "FIRE - AFTER WITHIN STATE ADJUSTMENT"
$ws.Range("E1:F1").Merge()
$ws.Range("C1:F1").HorizontalAlignment = -4108
$ws.Range("C1:F1").VerticalAlignment = -4108
$ws.Range("C1:F1").WrapText = $true
$ws.Range("A1").WrapText = $true
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").WrapText = $true

# border box around C1:D1 (matches "LIT RATES"/"NO FIRE" box style)
$box1 = $ws.Range("C1:D1")
$box1.Borders.Item(7).LineStyle = 1
$box1.Borders.Item(7).Weight = -4138
$box1.Borders.Item(10).LineStyle = 1
$box1.Borders.Item(10).Weight = -4138
$box1.Borders.Item(8).LineStyle = 1
$box1.Borders.Item(8).Weight = -4138
$box1.Borders.Item(9).LineStyle = 1
$box1.Borders.Item(9).Weight = 2
$ws.Range("C1").Borders.Item(10).LineStyle = 1
$ws.Range("C1").Borders.Item(10).Weight = -4138

$box2 = $ws.Range("E1:F1")
$box2.Borders.Item(7).LineStyle = 1
$box2.Borders.Item(7).Weight = -4138
$box2.Borders.Item(10).LineStyle = 1
$box2.Borders.Item(10).Weight = -4138
$box2.Borders.Item(8).LineStyle = 1
$box2.Borders.Item(8).Weight = -4138
$box2.Borders.Item(9).LineStyle = 1
$box2.Borders.Item(9).Weight = 2

# --- Header row 2 ---
$ws.Range("C2").Value = "% Change in HRU"
$ws.Range("D2").Value = "Mean proportions"
$ws.Range("E2").Value = "PROPORTIONS ON CYCLE 2"
$ws.Range("F2").Value = "Change from no fire"
$ws.Range("C2:F2").HorizontalAlignment = -4108
$ws.Range("C2:F2").VerticalAlignment = -4108
$ws.Range("C2:F2").WrapText = $true
$ws.Range("A2").WrapText = $true
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").WrapText = $true

$box3 = $ws.Range("C2:D2")
$box3.Borders.Item(7).LineStyle = 1
$box3.Borders.Item(7).Weight = -4138
$box3.Borders.Item(10).LineStyle = 1
$box3.Borders.Item(10).Weight = -4138
$box3.Borders.Item(9).LineStyle = 1
$box3.Borders.Item(9).Weight = -4138
$ws.Range("C2").Borders.Item(10).LineStyle = 1
$ws.Range("C2").Borders.Item(10).Weight = -4138

$box4 = $ws.Range("E2:F2")
$box4.Borders.Item(7).LineStyle = 1
$box4.Borders.Item(7).Weight = -4138
$box4.Borders.Item(10).LineStyle = 1
$box4.Borders.Item(10).Weight = -4138
$box4.Borders.Item(9).LineStyle = 1
$box4.Borders.Item(9).Weight = -4138

# --- Data rows 3-7 ---
$ws.Range("A3").Value = "seed=12345"
$ws.Range("B3").Value = "NONE"
$ws.Range("D3").Value = 0.92
$ws.Range("E3").Value = 0.8862
$ws.Range("F3").Formula = "=(E3-D3)/D3"

$ws.Range("B4").Value = "OCS"
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 0.038
$ws.Range("E4").Value = 0.053
$ws.Range("F4").Formula = "=(E4-D4)/D4"

$ws.Range("B5").Value = "UGT"
$ws.Range("D5").Value = 0.0338
$ws.Range("E5").Value = 0.044
$ws.Range("F5").Formula = "=(E5-D5)/D5"
$ws.Range("C4:C5").Merge()

$ws.Range("B6").Value = "ED"
$ws.Range("C6").Value = 1.12
$ws.Range("D6").Value = 0.0052
$ws.Range("E6").Value = 0.0114
$ws.Range("F6").Formula = "=(E6-D6)/D6"

$ws.Range("B7").Value = "HOSP"
$ws.Range("C7").Value = 0.53
$ws.Range("D7").Value = 0.0032
$ws.Range("E7").Value = 0.0054
$ws.Range("F7").Formula = "=(E7-D7)/D7"

# --- Formatting for data block B3:F7 ---
$ws.Range("B3:B7").HorizontalAlignment = -4108
$ws.Range("B3:B7").VerticalAlignment = -4108

$ws.Range("C3:C7").HorizontalAlignment = -4108
$ws.Range("C3:C7").VerticalAlignment = -4108

$dataBox = $ws.Range("B3:D7")
$dataBox.Borders.Item(7).LineStyle = 1
$dataBox.Borders.Item(7).Weight = -4138
$dataBox.Borders.Item(10).LineStyle = 1
$dataBox.Borders.Item(10).Weight = -4138
$ws.Range("B3:D3").Borders.Item(8).LineStyle = 1
$ws.Range("B3:D3").Borders.Item(8).Weight = -4138
$ws.Range("B7:D7").Borders.Item(9).LineStyle = 1
$ws.Range("B7:D7").Borders.Item(9).Weight = -4138
$ws.Range("C3:D7").Borders.Item(7).LineStyle = 1
$ws.Range("C3:D7").Borders.Item(7).Weight = -4138

$ws.Range("D3:D7").NumberFormat = "0.0000"
$ws.Range("D3:D7").HorizontalAlignment = -4108
$ws.Range("D3:D7").VerticalAlignment = -4108
for ($r = 3; $r -le 7; $r++) {
    $cell = $ws.Range("D$r")
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = -4138
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = -4138
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
}
$ws.Range("D7").Borders.Item(9).Weight = -4138

$ws.Range("E3:E6").NumberFormat = "0.0000"
$ws.Range("E3:E6").HorizontalAlignment = -4108
$ws.Range("E3:E6").VerticalAlignment = -4108
$ws.Range("E7").NumberFormat = "0.0000"
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").VerticalAlignment = -4108

$ws.Range("F3:F6").NumberFormat = "0%"
$ws.Range("F3:F6").HorizontalAlignment = -4108
$ws.Range("F3:F6").VerticalAlignment = -4108
$ws.Range("F3:F6").Font.Italic = $true
$ws.Range("F3:F6").Font.Color = 6787098

$ws.Range("F7").NumberFormat = "0%"
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F7").VerticalAlignment = -4108
$ws.Range("F7").Font.Italic = $true
$ws.Range("F7").Font.Color = 6787098

# --- Row 7 thick bottom border (closing border of the table) ---
$ws.Range("B7:F7").Borders.Item(9).LineStyle = 1
$ws.Range("B7:F7").Borders.Item(9).Weight = 4

$ws.Range("G11").Select() | Out-Null
